# Auto-generated script to apply Zodiark_Profits.xlsx market-data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 15325
$ws.Range("J58").Value = 17790
$ws.Range("L58").Value = 53370
$ws.Range("N58").Value = -53670
$ws.Range("H86").Value = 1639085.2
$ws.Range("I86").Value = 5503.091
$ws.Range("J86").Value = 4206143
$ws.Range("K86").Value = 5503.091
$ws.Range("L86").Value = 4206143
$ws.Range("M86").Value = -4380.091
$ws.Range("N86").Value = -4208389
$ws.Range("H89").Value = 1639085.2
$ws.Range("I89").Value = 5503.091
$ws.Range("J89").Value = 4206143
$ws.Range("K89").Value = 27515.455
$ws.Range("L89").Value = 21030715
$ws.Range("M89").Value = -21899.455
$ws.Range("N89").Value = -21041947
$ws.Range("H107").Value = 1376.1538
$ws.Range("I107").Value = 1060.2
$ws.Range("J107").Value = 2429.3333
$ws.Range("K107").Value = 1060.2
$ws.Range("L107").Value = 2429.3333
$ws.Range("M107").Value = 859.8
$ws.Range("N107").Value = -6269.3333
$ws.Range("H118").Value = 1028.1666
$ws.Range("I118").Value = 833.8
$ws.Range("K118").Value = 2501.4
$ws.Range("M118").Value = -844.3999999999996
$ws.Range("H129").Value = 1960.875
$ws.Range("I129").Value = 1883.8572
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 5651.571599999999
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = -651.5715999999993
$ws.Range("N129").Value = -17500
$ws.Range("H133").Value = 166743420
$ws.Range("J133").Value = 166743420
$ws.Range("L133").Value = 166743420
$ws.Range("N133").Value = -166753540
$ws.Range("H138").Value = 4903.859
$ws.Range("J138").Value = 5313.2026
$ws.Range("L138").Value = 15939.6078
$ws.Range("N138").Value = -26219.6078

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1629.7894
$ws.Range("I2").Value = 1570
$ws.Range("K2").Value = 1570
$ws.Range("M2").Value = -1457
$ws.Range("H32").Value = 3564.5686
$ws.Range("I32").Value = 2725
$ws.Range("J32").Value = 16997.666
$ws.Range("K32").Value = 2725
$ws.Range("L32").Value = 16997.666
$ws.Range("M32").Value = -2438
$ws.Range("N32").Value = -17571.666
$ws.Range("H43").Value = 24896.5
$ws.Range("J43").Value = 26530.334
$ws.Range("L43").Value = 26530.334
$ws.Range("N43").Value = -27156.334
$ws.Range("H45").Value = 1626.5834
$ws.Range("I45").Value = 1632.8572
$ws.Range("K45").Value = 1632.8572
$ws.Range("M45").Value = -1255.8572
$ws.Range("H61").Value = 3264.2666
$ws.Range("I61").Value = 3208.6316
$ws.Range("J61").Value = 3360.3635
$ws.Range("K61").Value = 3208.6316
$ws.Range("L61").Value = 3360.3635
$ws.Range("M61").Value = -2996.6316
$ws.Range("N61").Value = -3784.3635
$ws.Range("H63").Value = 2443.8
$ws.Range("I63").Value = 2443.8
$ws.Range("K63").Value = 2443.8
$ws.Range("M63").Value = -1757.8
$ws.Range("H66").Value = 2443.8
$ws.Range("I66").Value = 2443.8
$ws.Range("K66").Value = 12219
$ws.Range("M66").Value = -8787
$ws.Range("H116").Value = 1629.7894
$ws.Range("I116").Value = 1570
$ws.Range("K116").Value = 1570
$ws.Range("M116").Value = 724
$ws.Range("H132").Value = 7756.719
$ws.Range("I132").Value = 5424.5713
$ws.Range("J132").Value = 14286.733
$ws.Range("K132").Value = 16273.7139
$ws.Range("L132").Value = 42860.199
$ws.Range("M132").Value = -13743.7139
$ws.Range("N132").Value = -47920.199
$ws.Range("H134").Value = 137424.75
$ws.Range("J134").Value = 137424.75
$ws.Range("L134").Value = 137424.75
$ws.Range("N134").Value = -147564.75
$ws.Range("H136").Value = 3264.2666
$ws.Range("I136").Value = 3208.6316
$ws.Range("J136").Value = 3360.3635
$ws.Range("K136").Value = 9625.8948
$ws.Range("L136").Value = 10081.0905
$ws.Range("M136").Value = -7075.8948
$ws.Range("N136").Value = -15181.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1629.7894
$ws.Range("I3").Value = 1570
$ws.Range("K3").Value = 1570
$ws.Range("M3").Value = -1456
$ws.Range("H94").Value = 58826776
$ws.Range("I94").Value = 1739.2727
$ws.Range("J94").Value = 166672670
$ws.Range("K94").Value = 1739.2727
$ws.Range("L94").Value = 166672670
$ws.Range("M94").Value = -1288.2727
$ws.Range("N94").Value = -166673572
$ws.Range("H134").Value = 13515519
$ws.Range("I134").Value = 14707594
$ws.Range("K134").Value = 44122782
$ws.Range("M134").Value = -44120247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1256.5834
$ws.Range("I16").Value = 1235
$ws.Range("K16").Value = 1235
$ws.Range("M16").Value = -948
$ws.Range("H113").Value = 1256.5834
$ws.Range("I113").Value = 1235
$ws.Range("K113").Value = 1235
$ws.Range("M113").Value = 935
$ws.Range("H122").Value = 3198.4814
$ws.Range("I122").Value = 2847.4
$ws.Range("J122").Value = 3405
$ws.Range("K122").Value = 8542.200000000001
$ws.Range("L122").Value = 10215
$ws.Range("M122").Value = -6092.200000000001
$ws.Range("N122").Value = -15115
$ws.Range("H132").Value = 2302.375
$ws.Range("I132").Value = 2303.1428
$ws.Range("K132").Value = 6909.428400000001
$ws.Range("M132").Value = -4379.428400000001
$ws.Range("H134").Value = 1145.6842
$ws.Range("I134").Value = 1098.2778
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 3294.8334
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -759.8334000000004
$ws.Range("N134").Value = -11067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 562
$ws.Range("I17").Value = 179.4
$ws.Range("K17").Value = 538.2
$ws.Range("M17").Value = -369.2
$ws.Range("H34").Value = 1792
$ws.Range("J34").Value = 2199.4
$ws.Range("L34").Value = 6598.200000000001
$ws.Range("N34").Value = -6766.200000000001
$ws.Range("H39").Value = 4707.778
$ws.Range("J39").Value = 5208.75
$ws.Range("L39").Value = 15626.25
$ws.Range("N39").Value = -16214.25
$ws.Range("H55").Value = 2769.3
$ws.Range("I55").Value = 2833
$ws.Range("J55").Value = 2673.75
$ws.Range("K55").Value = 8499
$ws.Range("L55").Value = 8021.25
$ws.Range("M55").Value = -8322
$ws.Range("N55").Value = -8375.25
$ws.Range("H113").Value = 2931.4546
$ws.Range("J113").Value = 2594.8572
$ws.Range("L113").Value = 7784.571599999999
$ws.Range("N113").Value = -12124.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 41600
$ws.Range("I21").Value = 85000
$ws.Range("K21").Value = 85000
$ws.Range("M21").Value = -84827
$ws.Range("H30").Value = 41600
$ws.Range("I30").Value = 85000
$ws.Range("K30").Value = 85000
$ws.Range("M30").Value = -84895
$ws.Range("H33").Value = 34000
$ws.Range("J33").Value = 34000
$ws.Range("L33").Value = 34000
$ws.Range("N33").Value = -34504
$ws.Range("H70").Value = 44735.188
$ws.Range("I70").Value = 82929.86
$ws.Range("J70").Value = 15028.223
$ws.Range("K70").Value = 82929.86
$ws.Range("L70").Value = 15028.223
$ws.Range("M70").Value = -82659.86
$ws.Range("N70").Value = -15568.223
$ws.Range("H73").Value = 44735.188
$ws.Range("I73").Value = 82929.86
$ws.Range("J73").Value = 15028.223
$ws.Range("K73").Value = 82929.86
$ws.Range("L73").Value = 15028.223
$ws.Range("M73").Value = -81993.86
$ws.Range("N73").Value = -16900.223
$ws.Range("H102").Value = 2975
$ws.Range("I102").Value = 2770.2
$ws.Range("K102").Value = 2770.2
$ws.Range("M102").Value = -1148.2
$ws.Range("H126").Value = 8121
$ws.Range("I126").Value = 7535
$ws.Range("K126").Value = 22605
$ws.Range("M126").Value = -20135

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1969.1177
$ws.Range("I61").Value = 2001.0714
$ws.Range("K61").Value = 2001.0714
$ws.Range("M61").Value = -1799.0714
$ws.Range("H68").Value = 4954.778
$ws.Range("I68").Value = 1527.5
$ws.Range("J68").Value = 11809.333
$ws.Range("K68").Value = 1527.5
$ws.Range("L68").Value = 11809.333
$ws.Range("M68").Value = -778.5
$ws.Range("N68").Value = -13307.333
$ws.Range("H71").Value = 4954.778
$ws.Range("I71").Value = 1527.5
$ws.Range("J71").Value = 11809.333
$ws.Range("K71").Value = 7637.5
$ws.Range("L71").Value = 59046.665
$ws.Range("M71").Value = -3893.5
$ws.Range("N71").Value = -66534.66500000001
$ws.Range("H113").Value = 1969.1177
$ws.Range("I113").Value = 2001.0714
$ws.Range("K113").Value = 2001.0714
$ws.Range("M113").Value = 168.9286
$ws.Range("H136").Value = 3498.5908
$ws.Range("I136").Value = 3161.762
$ws.Range("J136").Value = 3806.1304
$ws.Range("K136").Value = 9485.286
$ws.Range("L136").Value = 11418.3912
$ws.Range("M136").Value = -6935.286
$ws.Range("N136").Value = -16518.3912

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10438
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10438
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 20876
$ws.Range("N81").Value = -22998
$ws.Range("H84").Value = 10438
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10438
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 104380
$ws.Range("N84").Value = -114988
$ws.Range("H97").Value = 12133.333
$ws.Range("J97").Value = 12133.333
$ws.Range("L97").Value = 12133.333
$ws.Range("N97").Value = -14115.333
$ws.Range("H107").Value = 804.8
$ws.Range("I107").Value = 668.5
$ws.Range("J107").Value = 1009.25
$ws.Range("K107").Value = 2005.5
$ws.Range("L107").Value = 3027.75
$ws.Range("M107").Value = -85.5
$ws.Range("N107").Value = -6867.75
$ws.Range("H132").Value = 1098.7333
$ws.Range("I132").Value = 694.6
$ws.Range("K132").Value = 2083.8
$ws.Range("M132").Value = 446.1999999999998
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()
